# Finish PostProcessing with VBA
#
# Corrects the Factory/Product numbers that were mixed up between the
# "Hai Phong" and "Binh Dinh" factories on the Inbound sheets, and adds the
# previously-missing "Binh Dinh" / "Bulk" row that PostProcessing produces.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Inbound Cost Per Factory")
$ws2 = $wb.Worksheets.Item("Outbound Cost Per Customer")
$ws3 = $wb.Worksheets.Item("Inbound Volume Per Factory")
$ws4 = $wb.Worksheets.Item("Outbound Volume Per Customer")

# --- "Inbound Cost Per Factory" ---
$ws1.Range("A2").Value = "Binh Dinh"
$ws1.Range("B2").Value = "Bag"
$ws1.Range("C2").Value = 600
$ws1.Range("D2").Value = 1260

$ws1.Range("A3").Value = "Hai Phong"
$ws1.Range("B3").Value = "Bag"
$ws1.Range("C3").Value = 300
$ws1.Range("D3").Value = 360

$ws1.Range("A4").Value = "Binh Dinh"
$ws1.Range("B4").Value = "Bulk"
$ws1.Range("C4").Value = 1350
$ws1.Range("D4").Value = 1530

# --- "Inbound Volume Per Factory" ---
$ws3.Range("A2").Value = "Binh Dinh"
$ws3.Range("B2").Value = "Bag"
$ws3.Range("C2").Value = 20
$ws3.Range("D2").Value = 42

$ws3.Range("A3").Value = "Hai Phong"
$ws3.Range("B3").Value = "Bag"
$ws3.Range("C3").Value = 10
$ws3.Range("D3").Value = 12

$ws3.Range("A4").Value = "Binh Dinh"
$ws3.Range("B4").Value = "Bulk"
$ws3.Range("C4").Value = 45
$ws3.Range("D4").Value = 51

# --- Reset view/selection state on the sheets that are no longer the
# --- active tab, then make "Inbound Cost Per Factory" the active sheet
# --- with cell H16 selected, matching the finished workbook.
[void]$ws3.Range("A1").Select()
[void]$ws4.Range("A1").Select()

[void]$ws1.Select()
[void]$ws1.Range("H16").Select()
